$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '30.418.11'
$ws.Range('E2').Value = '  -0.39%  '
Set-TextValue 'D3' '2.099.43'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue 'D5' '333.95'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('E6').Value = '  +0.15%  '
Set-TextValue 'D7' '0.5204'
$ws.Range('E7').Value = '  -1.18%  '
Set-TextValue 'D8' '0.4535'
$ws.Range('E8').Value = '  +3.84%  '
Set-TextValue 'D9' '54.35'
$ws.Range('E9').Value = '  +14.72%  '
Set-TextValue 'D10' '0.08879'
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('E11').Value = '  +1.30%  '
Set-TextValue 'D12' '24.06'
$ws.Range('E12').Value = '  -2.05%  '
Set-TextValue 'D13' '2.095.39'
$ws.Range('E13').Value = '  -0.55%  '
Set-TextValue 'D14' '6.795'
$ws.Range('E14').Value = '  +0.95%  '
Set-TextValue 'D15' '8.018'
$ws.Range('E15').Value = '  +3.22%  '
Set-TextValue 'D16' '96.84'
$ws.Range('E16').Value = '  +0.47%  '
Set-TextValue 'D17' '0.00001141'
$ws.Range('E17').Value = '  +1.10%  '
Set-TextValue 'D18' '1.003'
$ws.Range('E18').Value = '  +0.07%  '
Set-TextValue 'D19' '0.06618'
$ws.Range('E19').Value = '  -0.32%  '
Set-TextValue 'D20' '19.16'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  +0.18%  '
Set-TextValue 'D22' '6.268'
$ws.Range('E22').Value = '  -0.87%  '
Set-TextValue 'D23' '30.454.66'
$ws.Range('E23').Value = '  -0.44%  '
Set-TextValue 'D24' '12.32'
$ws.Range('E24').Value = '  +0.30%  '
Set-TextValue 'D25' '2.334'
$ws.Range('E25').Value = '  -1.02%  '
Set-TextValue 'D26' '2.338.56'
$ws.Range('E26').Value = '  -0.65%  '
Set-TextValue 'D27' '22.16'
$ws.Range('E27').Value = '  -1.20%  '
Set-TextValue 'D28' '162.82'
$ws.Range('E28').Value = '  +0.51%  '
Set-TextValue 'D29' '2.512'
$ws.Range('E29').Value = '  -3.25%  '
Set-TextValue 'D30' '132.91'
$ws.Range('E30').Value = '  +0.16%  '
Set-TextValue 'D31' '1.202'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -1.07%  '
Set-TextValue 'D33' '1.653'
$ws.Range('E33').Value = '  -1.38%  '
Set-TextValue 'D34' '6.361'
Set-TextValue 'D35' '3.942'
$ws.Range('E35').Value = '  +0.38%  '
Set-TextValue 'D36' '10.39'
$ws.Range('E36').Value = '  +3.39%  '
Set-TextValue 'D37' '5.795'
$ws.Range('E37').Value = '  +5.42%  '
Set-TextValue 'D38' '0.02570'
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.2296'
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D41' '12.72'
$ws.Range('E41').Value = '  -0.32%  '
Set-TextValue 'D42' '0.6861'
$ws.Range('E42').Value = '  +0.78%  '
Set-TextValue 'D43' '1.246'
$ws.Range('E43').Value = '  -1.17%  '
Set-TextValue 'D44' '2.318'
$ws.Range('E44').Value = '  +4.90%  '
Set-TextValue 'D45' '13.92'
$ws.Range('E45').Value = '  -0.69%  '
Set-TextValue 'D46' '0.6339'
$ws.Range('E46').Value = '  -0.45%  '
Set-TextValue 'D47' '3.657'
$ws.Range('E47').Value = '  +0.86%  '
Set-TextValue 'D48' '1.245'
$ws.Range('E48').Value = '  -0.53%  '
Set-TextValue 'D49' '0.00000000346'
$ws.Range('E49').Value = '  +17.61%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '82.90'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D51' '1.199'
$ws.Range('E51').Value = '  +0.16%  '
